{"js": "const body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst replacements = [\n  { row: 0, col: 0, oldText: \"95\u00f72=47, 1\", newText: \"72\u00f79=8, 0\" },\n  { row: 0, col: 1, oldText: \"10\u00f72=5, 0\", newText: \"10\u00f78=1, 2\" },\n  { row: 0, col: 2, oldText: \"61\u00f78=7, 5\", newText: \"26\u00f78=3, 2\" },\n  { row: 0, col: 3, oldText: \"89\u00f73=29, 2\", newText: \"32\u00f72=16, 0\" },\n  { row: 0, col: 4, oldText: \"92\u00f77=13, 1\", newText: \"58\u00f76=9, 4\" },\n  { row: 4, col: 0, oldText: \"26\u00f75=5, 1\", newText: \"41\u00f75=8, 1\" },\n  { row: 4, col: 1, oldText: \"39\u00f72=19, 1\", newText: \"61\u00f75=12, 1\" },\n  { row: 4, col: 2, oldText: \"68\u00f77=9, 5\", newText: \"99\u00f74=24, 3\" },\n  { row: 4, col: 3, oldText: \"64\u00f74=16, 0\", newText: \"84\u00f74=21, 0\" },\n  { row: 4, col: 4, oldText: \"98\u00f73=32, 2\", newText: \"86\u00f79=9, 5\" },\n  { row: 8, col: 0, oldText: \"61\u00f73=20, 1\", newText: \"35\u00f74=8, 3\" },\n  { row: 8, col: 1, oldText: \"66\u00f79=7, 3\", newText: \"76\u00f72=38, 0\" },\n  { row: 8, col: 2, oldText: \"39\u00f76=6, 3\", newText: \"34\u00f75=6, 4\" },\n  { row: 8, col: 3, oldText: \"42\u00f74=10, 2\", newText: \"58\u00f74=14, 2\" },\n  { row: 8, col: 4, oldText: \"99\u00f75=19, 4\", newText: \"73\u00f76=12, 1\" },\n  { row: 12, col: 0, oldText: \"94\u00f73=31, 1\", newText: \"85\u00f76=14, 1\" },\n  { row: 12, col: 1, oldText: \"59\u00f72=29, 1\", newText: \"98\u00f73=32, 2\" },\n  { row: 12, col: 2, oldText: \"81\u00f75=16, 1\", newText: \"95\u00f77=13, 4\" },\n  { row: 12, col: 3, oldText: \"27\u00f78=3, 3\", newText: \"23\u00f76=3, 5\" },\n  { row: 12, col: 4, oldText: \"58\u00f75=11, 3\", newText: \"30\u00f73=10, 0\" },\n  { row: 16, col: 0, oldText: \"14\u00f79=1, 5\", newText: \"16\u00f76=2, 4\" },\n  { row: 16, col: 1, oldText: \"61\u00f77=8, 5\", newText: \"62\u00f79=6, 8\" },\n  { row: 16, col: 2, oldText: \"68\u00f77=9, 5\", newText: \"83\u00f75=16, 3\" },\n  { row: 16, col: 3, oldText: \"78\u00f79=8, 6\", newText: \"88\u00f72=44, 0\" },\n  { row: 16, col: 4, oldText: \"80\u00f72=40, 0\", newText: \"85\u00f77=12, 1\" },\n];\n\nfor (const r of replacements) {\n  const cell = table.getCell(r.row, r.col);\n  cell.body.load(\"text\");\n  await context.sync();\n  if (cell.body.text !== r.oldText) {\n    throw new Error(\n      \"Cell (\" + r.row + \",\" + r.col + \") expected '\" + r.oldText + \"' but found '\" + cell.body.text + \"'\"\n    );\n  }\n\n  const results = cell.body.search(r.oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for \" + r.oldText + \" at row \" + r.row + \" col \" + r.col);\n  }\n  results.items[0].insertText(r.newText, \"Replace\");\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n    @{ Row = 1; Col = 1; OldText = \"95\u00f72=47, 1\"; NewText = \"72\u00f79=8, 0\" },\n    @{ Row = 1; Col = 2; OldText = \"10\u00f72=5, 0\"; NewText = \"10\u00f78=1, 2\" },\n    @{ Row = 1; Col = 3; OldText = \"61\u00f78=7, 5\"; NewText = \"26\u00f78=3, 2\" },\n    @{ Row = 1; Col = 4; OldText = \"89\u00f73=29, 2\"; NewText = \"32\u00f72=16, 0\" },\n    @{ Row = 1; Col = 5; OldText = \"92\u00f77=13, 1\"; NewText = \"58\u00f76=9, 4\" },\n    @{ Row = 5; Col = 1; OldText = \"26\u00f75=5, 1\"; NewText = \"41\u00f75=8, 1\" },\n    @{ Row = 5; Col = 2; OldText = \"39\u00f72=19, 1\"; NewText = \"61\u00f75=12, 1\" },\n    @{ Row = 5; Col = 3; OldText = \"68\u00f77=9, 5\"; NewText = \"99\u00f74=24, 3\" },\n    @{ Row = 5; Col = 4; OldText = \"64\u00f74=16, 0\"; NewText = \"84\u00f74=21, 0\" },\n    @{ Row = 5; Col = 5; OldText = \"98\u00f73=32, 2\"; NewText = \"86\u00f79=9, 5\" },\n    @{ Row = 9; Col = 1; OldText = \"61\u00f73=20, 1\"; NewText = \"35\u00f74=8, 3\" },\n    @{ Row = 9; Col = 2; OldText = \"66\u00f79=7, 3\"; NewText = \"76\u00f72=38, 0\" },\n    @{ Row = 9; Col = 3; OldText = \"39\u00f76=6, 3\"; NewText = \"34\u00f75=6, 4\" },\n    @{ Row = 9; Col = 4; OldText = \"42\u00f74=10, 2\"; NewText = \"58\u00f74=14, 2\" },\n    @{ Row = 9; Col = 5; OldText = \"99\u00f75=19, 4\"; NewText = \"73\u00f76=12, 1\" },\n    @{ Row = 13; Col = 1; OldText = \"94\u00f73=31, 1\"; NewText = \"85\u00f76=14, 1\" },\n    @{ Row = 13; Col = 2; OldText = \"59\u00f72=29, 1\"; NewText = \"98\u00f73=32, 2\" },\n    @{ Row = 13; Col = 3; OldText = \"81\u00f75=16, 1\"; NewText = \"95\u00f77=13, 4\" },\n    @{ Row = 13; Col = 4; OldText = \"27\u00f78=3, 3\"; NewText = \"23\u00f76=3, 5\" },\n    @{ Row = 13; Col = 5; OldText = \"58\u00f75=11, 3\"; NewText = \"30\u00f73=10, 0\" },\n    @{ Row = 17; Col = 1; OldText = \"14\u00f79=1, 5\"; NewText = \"16\u00f76=2, 4\" },\n    @{ Row = 17; Col = 2; OldText = \"61\u00f77=8, 5\"; NewText = \"62\u00f79=6, 8\" },\n    @{ Row = 17; Col = 3; OldText = \"68\u00f77=9, 5\"; NewText = \"83\u00f75=16, 3\" },\n    @{ Row = 17; Col = 4; OldText = \"78\u00f79=8, 6\"; NewText = \"88\u00f72=44, 0\" },\n    @{ Row = 17; Col = 5; OldText = \"80\u00f72=40, 0\"; NewText = \"85\u00f77=12, 1\" },\n)\n\nforeach ($r in $replacements) {\n    $cell = $t.Cell($r.Row, $r.Col)\n    $rng = $cell.Range\n    $current = $rng.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $r.OldText) {\n        throw \"Cell ($($r.Row),$($r.Col)) expected `\"$($r.OldText)`\" but found `\"$current`\"\"\n    }\n    $rng.Text = $r.NewText\n}"}
